# Update omics and assembly tags on the MIMS worksheet (isa_template / sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a second TAGS entry ("Assembly") alongside the existing "Metagenomics" tag,
# together with its term accession number and term source REF.
$ws.Range("E13").Value = "Assembly"
$ws.Range("E14").Value = "https://bioregistry.io/NCIT:C52474"
$ws.Range("E15").Value = "NCIT"

# Reflect the new selection state left behind by the edit (select the new cells).
$ws.Range("E13:E15").Select()
